# BAU Zero Emis Credit for Nuc Plants.xlsx
# Commit: "test elec sec US files"
#
# The authoring change strips the per-state lookup helper data out of the
# "About" sheet (the U.S. state name/abbreviation table that lived in
# columns F:G plus the supporting B1/C1/B2 helper cells), and removes the
# now-unused "Sheet1" worksheet that held a small scratch copy of that
# lookup table. The remaining "BZECfNP" worksheet becomes the second (and
# now last) tab, and the "About" tab becomes the active tab/selection.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$about = $wb.Worksheets.Item("About")
$bz    = $wb.Worksheets.Item("BZECfNP")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- Clean out the helper / lookup-table content on "About" ---
# F1:G50 held the full state name -> abbreviation lookup table
$about.Range("F1:G50").Clear()
# B1 (selected state name) and C1 (date stamp) were inputs tied to that table
$about.Range("B1:C1").Clear()
# B2 held =LOOKUP(B1,F1:G50,G1:G50) resolving the state abbreviation
$about.Range("B2").Clear()

# --- Remove the now unused "Sheet1" scratch worksheet ---
$sheet1.Delete()

# --- Update view/selection state to match the edited workbook ---
# BZECfNP: scroll/selection moved down to the last data row
$bz.Activate()
$bz.Range("B25:AE25").Select()

# "About" becomes the active/selected tab of the workbook
$about.Activate()
$about.Range("A1").Select()
